# Add columns I (I0) and J (IF) to the sheet, mirroring the header/style
# treatment of the existing columns (B1:H1 use the bold/bordered/centered
# style), then fill in the per-row numeric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header style used by the other header cells (B1:H1) by copying
# H1's format onto the new header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for I2:J11
$iValues = @(6, 9, 4, 7, 9, 9, 5, 9, 3, 7)
$jValues = @(6, 9, 5, 7, 9, 9, 5, 9, 3, 7)

for ($r = 0; $r -lt 10; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$r]
    $ws.Cells.Item($row, 10).Value = $jValues[$r]
}
